$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.870.65"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "'3.040.41"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'580.27"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").Value = "'150.35"
$ws.Range("E6").Value = "  -3.07%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.528"
$ws.Range("E8").Value = "  -2.52%  "
$ws.Range("D9").Value = "'3.042.02"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("E10").Value = "  -3.40%  "
$ws.Range("D11").Value = "'5.76"
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").Value = "'0.445"
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("E13").Value = "  -3.71%  "
$ws.Range("D14").Value = "'35.72"
$ws.Range("E14").Value = "  -4.73%  "
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "'3.542.50"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("D17").Value = "'7.07"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").Value = "'62.829.13"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").Value = "'3.040.28"
$ws.Range("E19").Value = "  -1.46%  "
$ws.Range("D20").Value = "'475.55"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "'14.19"
$ws.Range("E21").Value = "  -3.00%  "
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("D23").Value = "'7.45"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").Value = "'2.35"
$ws.Range("E24").Value = "  -2.57%  "
$ws.Range("D25").Value = "'81.16"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'12.53"
$ws.Range("E26").Value = "  -3.07%  "
$ws.Range("D27").Value = "'10.50"
$ws.Range("E27").Value = "  +4.26%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.26"
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("D32").Value = "'2.18"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").Value = "'27.50"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").Value = "'0.0₃0801"
$ws.Range("E36").Value = "  -5.89%  "
$ws.Range("D37").Value = "'5.83"
$ws.Range("E37").Value = "  -4.23%  "
$ws.Range("E38").Value = "  -2.92%  "
$ws.Range("D39").Value = "'3.07"
$ws.Range("E39").Value = "  -9.41%  "
$ws.Range("D40").Value = "'50.14"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("D41").Value = "'9.09"
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("D42").Value = "'422.65"
$ws.Range("E42").Value = "  -4.44%  "
$ws.Range("D43").Value = "'0.283"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").Value = "'0.113"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("D45").Value = "'2.823.67"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("D46").Value = "'0.0358"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").Value = "'37.90"
$ws.Range("E47").Value = "  -6.09%  "
$ws.Range("D48").Value = "'127.05"
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("D50").Value = "'24.75"
$ws.Range("E50").Value = "  -3.90%  "
$ws.Range("E51").Value = "  -1.15%  "
